# 1401ME09 marksheet — recompute the score summary (rows 10-12) and populate
# the "Student Ans" column with the option the student actually picked
# (rows 16-40), then drop the now-unused 3rd Student/Correct-Ans block
# (columns G:H) together with the 2nd block (D:E) for every question row
# that doesn't need it any more.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Score summary block
# ---------------------------------------------------------------------

# Row 10 ("No.") / Row 11 ("Marking") / Row 12 ("Total") labels pick up the
# same centred "mtitleStyle" formatting already used by the header row (A9).
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B10").Value = 12
$ws.Range("D10").Value = 16
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 48
$ws.Range("E12").Value = "48/112"

# ---------------------------------------------------------------------
# Per-question "Student Ans" column: fill in column A (and D for row 18)
# with the option the student chose, using the same "correctStyle"
# highlight already applied to the matching "Correct Ans" cell.
# ---------------------------------------------------------------------

$studentPicks = @{
    16 = "Option A"
    18 = "Option B"
    19 = "Option C"
    20 = "Option B"
    21 = "Option C"
    24 = "Option A"
    25 = "Option A"
    29 = "Option D"
    32 = "Option C"
    33 = "Option D"
    35 = "Option D"
}

foreach ($row in $studentPicks.Keys) {
    $ws.Range("B$row").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").Value = $studentPicks[$row]
}

# Row 18 also keeps its 2nd block, so mirror the same pick into D18.
$ws.Range("E18").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Option D"

# ---------------------------------------------------------------------
# Drop the unused 2nd block (D:E) for every row that no longer needs it,
# and the unused 3rd block (G:H) entirely.
# ---------------------------------------------------------------------

$ws.Range("D19:E40").Clear()
$ws.Range("G15:H21").Clear()
